## Append new scrape result: 2025-11-20 18:26 JST
## - A new top-of-list job ("【地域活性化】自社サイトにソーシャルウォールを導入したい")
##   was found, so it is inserted at row 13 and everything that used to be at
##   row 13 onward shifts down by one row.
## - Every data row's "取得日時" (fetched-at) timestamp is refreshed to the
##   new scrape time.
## - The worksheet's hyperlinks (URL column F) have to be rebuilt because the
##   row insert does not re-point the existing hyperlink relationships.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-20 18:26:42"

# 1) Insert a new row at position 13. Excel shifts rows 13-22 down to 14-23,
#    carrying their values/styles/formatting with them.
$ws.Range("A13").EntireRow.Insert()

# 2) Populate the newly inserted row 13 with the freshly scraped listing.
$ws.Cells.Item(13, 1).Value = $newTimestamp
$ws.Cells.Item(13, 2).Value = "【地域活性化】自社サイトにソーシャルウォールを導入したい"
$ws.Cells.Item(13, 3).Value = "システム開発"
$ws.Cells.Item(13, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(13, 5).Value = "期限情報なし"
$ws.Cells.Item(13, 6).Value = "https://www.lancers.jp/work/detail/5438358"
$ws.Cells.Item(13, 7).Value = 33
$ws.Cells.Item(13, 8).Value = "◇サイト"

# 3) Refresh the "取得日時" timestamp on every data row (2-23) to the new
#    scrape time (the rest of rows 2-12 and 14-23 keep their prior content).
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# 4) Rebuild the hyperlinks for column F (rows 2-23). The row-insert above
#    leaves the sheet's hyperlink relationships pointing at the pre-insert
#    rows, so clear them all out and re-add them against the final layout.
$ws.Range("A1").Hyperlinks.Delete()

$urls = @(
    "https://www.lancers.jp/work/detail/5437916",
    "https://www.lancers.jp/work/detail/5437447",
    "https://www.lancers.jp/work/detail/5437717",
    "https://www.lancers.jp/work/detail/5427956",
    "https://www.lancers.jp/work/detail/5437726",
    "https://www.lancers.jp/work/detail/5438035",
    "https://www.lancers.jp/work/detail/5438164",
    "https://www.lancers.jp/work/detail/5438081",
    "https://www.lancers.jp/work/detail/5437832",
    "https://www.lancers.jp/work/detail/5437868",
    "https://www.lancers.jp/work/detail/5437655",
    "https://www.lancers.jp/work/detail/5438358",
    "https://www.lancers.jp/work/detail/5437728",
    "https://www.lancers.jp/work/detail/5437997",
    "https://www.lancers.jp/work/detail/5437991",
    "https://www.lancers.jp/work/detail/5429335",
    "https://www.lancers.jp/work/detail/5437783",
    "https://www.lancers.jp/work/detail/5437544",
    "https://www.lancers.jp/work/detail/5432661",
    "https://www.lancers.jp/work/detail/5438052",
    "https://www.lancers.jp/work/detail/5436248",
    "https://www.lancers.jp/work/detail/5438014"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $urls[$i]
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $urls[$i])
}
